# Threshold/Zn/3His/4.xlsx — apply the recorded edit:
#   - alpha_distance_range Min  5.6 -> 4
#   - beta_distance_range  Min  5.7 -> 5
#   - ratio_threshold_range Min 0.9 -> 0.7, Max 1.2 -> 1.3
#   - theta_threshold_range row removed entirely (its shared string is dropped too)
#   - pie_threshold_range shifts up from row 6 to row 5 (values unchanged: 0 / 20)
#   - sheet view selection now rests on C5
#   - a page setup (Letter portrait) is recorded for the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the changed numeric values in place -----------------------------
$ws.Cells.Item(2, 2).Value = 4      # B2  alpha_distance_range / Min
$ws.Cells.Item(3, 2).Value = 5      # B3  beta_distance_range  / Min
$ws.Cells.Item(4, 2).Value = 0.7    # B4  ratio_threshold_range / Min
$ws.Cells.Item(4, 3).Value = 1.3    # C4  ratio_threshold_range / Max

# --- Remove the theta_threshold_range row (row 5) ----------------------------
# This shifts the former row 6 (pie_threshold_range, Min 0 / Max 20) up to
# row 5 and drops "theta_threshold_range" out of the shared-string table.
$ws.Rows(5).Delete()

# --- Page setup recorded for the sheet --------------------------------------
$ws.PageSetup.PaperSize = 9       # xlPaperA4 (paperSize="9")
$ws.PageSetup.Orientation = 1     # xlPortrait

# --- Selection state ----------------------------------------------------------
$ws.Range("C5").Select()
